# Generate Report for Handoff
#
# Refreshes the localization-status report after a new handoff report was
# generated:
#   - The "Priority" for the e2e markdown source rows (4-7) changes from
#     "low" to "ht" on both the "zh-cn" and "de-de" sheets.
#   - The recorded handoff timestamps for those same rows are refreshed to
#     reflect the new report generation time.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newHoGenerateDate  = "2016-08-29 20:43:59"
$newZhHandoffDate   = "2016-08-29 20:43:54"

foreach ($r in 4..7) {
    # Overview: "Latest HO Xliff Generate Date" column (G)
    $wsOverview.Range("G$r").Value = $newHoGenerateDate

    # zh-cn: "Priority" column (E) and "Latest Handoff Datetime" column (H)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = $newZhHandoffDate

    # de-de: "Priority" column (E) and "Latest Handoff Datetime" column (H)
    # (de-de's handoff datetime shares the new Overview report timestamp)
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = $newHoGenerateDate
}
